$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns (D, E) are stored as literal text in the source data
# (e.g. "30.521.24", "0.08760", "  +1.71%  ") even when numeric-looking.
# A bare .Value assignment lets Excel auto-detect such strings as Numbers
# and normalize them (stripping significant trailing zeros, losing the
# original formatting). Force Text format before assigning, then restore
# the default "Normal" cell style afterwards so no stray style survives.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "D29",
    "E29",
    "D30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "D35",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51",
)
foreach ($cell in $textCells) { $ws.Range($cell).NumberFormat = "@" }

$ws.Range("D2").Value = '30.521.24'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").Value = '2.018.92'
$ws.Range("E3").Value = '  +5.75%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '324.78'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.5125'
$ws.Range("E7").Value = '  +1.60%  '
$ws.Range("D8").Value = '0.4236'
$ws.Range("E8").Value = '  +4.52%  '
$ws.Range("D9").Value = '0.08760'
$ws.Range("E9").Value = '  +5.81%  '
$ws.Range("D10").Value = '1.138'
$ws.Range("E10").Value = '  +3.50%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '43.40'
$ws.Range("E11").Value = '  +3.30%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '24.76'
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '2.015.30'
$ws.Range("E13").Value = '  +5.59%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '6.617'
$ws.Range("E14").Value = '  +3.64%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.493'
$ws.Range("E15").Value = '  +3.38%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '94.51'
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.00001115'
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.06534'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '19.03'
$ws.Range("E20").Value = '  +5.43%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.232'
$ws.Range("E22").Value = '  +5.03%  '
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '30.589.85'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '11.88'
$ws.Range("E24").Value = '  +5.32%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.231'
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.251.51'
$ws.Range("E26").Value = '  +5.72%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '22.45'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '163.13'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.439'
$ws.Range("E29").Value = '  +7.23%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '131.65'
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.143'
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.1053'
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.099'
$ws.Range("E33").Value = '  +2.56%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.833'
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.365'
$ws.Range("E35").Value = '  +14.52%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02528'
$ws.Range("E36").Value = '  +3.78%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.509'
$ws.Range("E37").Value = '  +2.04%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.06656'
$ws.Range("E38").Value = '  +5.00%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '12.30'
$ws.Range("E39").Value = '  +8.41%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '9.117'
$ws.Range("E40").Value = '  +5.05%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2202'
$ws.Range("E41").Value = '  +2.42%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6694'
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.232'
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.69'
$ws.Range("E45").Value = '  +3.27%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6199'
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.203'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.661'
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.270'
$ws.Range("E49").Value = '  +5.10%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '124.88'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '81.25'
$ws.Range("E51").Value = '  +3.81%  '

foreach ($cell in $textCells) { $ws.Range($cell).Style = "Normal" }
